# Updated Exercises 2 results 23.11.2016
# Fills in column D ("Оценка котролно 2") and column E ("Обща оценка")
# for rows 2-19, using "N/A" for the two students who did not take the
# second exercise (rows 6 and 15). Those "N/A" cells are right-aligned.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> (D value, E value)
$values = @{
    2  = @(6, 6)
    3  = @(6, 6)
    4  = @(6, 6)
    5  = @(6, 6)
    6  = @("N/A", 4)
    7  = @(6, 4)
    8  = @(6, 6)
    9  = @(6, 6)
    10 = @(6, 4)
    11 = @(6, 5)
    12 = @(5, 5)
    13 = @(5, 5)
    14 = @(6, 5)
    15 = @("N/A", 4)
    16 = @(6, 6)
    17 = @(5, 4)
    18 = @(5, 4)
    19 = @(6, 6)
}

foreach ($row in 2..19) {
    $pair = $values[$row]
    $dCell = $ws.Cells.Item($row, 4)
    $dCell.Value = $pair[0]
    if ($pair[0] -eq "N/A") {
        $dCell.HorizontalAlignment = -4152
    }
    $ws.Cells.Item($row, 5).Value = $pair[1]
}

# Move the active selection to B1 (was E22)
$null = $ws.Range("B1").Select()
